# Nuevo formato 15 jun 2021
# Update individual grade cells in the "Calificaciones" sheet that were
# previously placeholders (-1) with the real computed grade values, plus a
# handful of other grade corrections, as captured in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calificaciones")

$ws.Range("N4").Value = 7
$ws.Range("N5").Value = 8
$ws.Range("N6").Value = 7
$ws.Range("T6").Value = 7
$ws.Range("N7").Value = 5
$ws.Range("N8").Value = 7
$ws.Range("N9").Value = 6
$ws.Range("N10").Value = 7
$ws.Range("N11").Value = 9
$ws.Range("N12").Value = 5
$ws.Range("N13").Value = 7
$ws.Range("N14").Value = 5
$ws.Range("N15").Value = 9
$ws.Range("N16").Value = 8
$ws.Range("N17").Value = 7
$ws.Range("N18").Value = 8
$ws.Range("T18").Value = 7
$ws.Range("N19").Value = 7
$ws.Range("N20").Value = 9
$ws.Range("T20").Value = 8
$ws.Range("N21").Value = 8
$ws.Range("S21").Value = 9
$ws.Range("T21").Value = 6
$ws.Range("Y21").Value = 6
$ws.Range("N22").Value = 7
$ws.Range("T22").Value = 6
$ws.Range("N23").Value = 7
$ws.Range("N24").Value = 7
$ws.Range("S24").Value = 8
$ws.Range("T24").Value = 6
$ws.Range("Y24").Value = 6
$ws.Range("N25").Value = 6
$ws.Range("N26").Value = 5
$ws.Range("N27").Value = 7
$ws.Range("T27").Value = 6
$ws.Range("N28").Value = 8
$ws.Range("T28").Value = 8
$ws.Range("N29").Value = 6
$ws.Range("S29").Value = 5
$ws.Range("Y29").Value = 5
$ws.Range("N30").Value = 7
$ws.Range("S30").Value = 9
$ws.Range("T30").Value = 6
$ws.Range("Y30").Value = 6
$ws.Range("N31").Value = 7
$ws.Range("T31").Value = 8
$ws.Range("N32").Value = 7
$ws.Range("N33").Value = 6
$ws.Range("N34").Value = 6
$ws.Range("N35").Value = 7
$ws.Range("N36").Value = 7
$ws.Range("T36").Value = 6
$ws.Range("N37").Value = 6
$ws.Range("N38").Value = 9
$ws.Range("N39").Value = 8
